{"js": "// Replace the date line and each \"A\u00d7B=C\" table-cell expression with its\n// updated value, per the commit diff (26 text substitutions total: 1 date\n// heading + 25 multiplication-table cells). Each \"old\" string is a unique\n// substring of the document, so body.search(old, {matchCase:true}) always\n// returns exactly one hit; insertText(new, \"Replace\") swaps the run's text\n// in place without touching its formatting (font/size/alignment).\nconst replacements = [\n  [\"2025-09-22 Monday\", \"2025-09-23 Tuesday\"],\n  [\"516\u00d77=3612\", \"716\u00d73=2148\"],\n  [\"425\u00d72=850\", \"467\u00d75=2335\"],\n  [\"838\u00d73=2514\", \"135\u00d73=405\"],\n  [\"824\u00d74=3296\", \"715\u00d73=2145\"],\n  [\"531\u00d77=3717\", \"245\u00d79=2205\"],\n  [\"154\u00d77=1078\", \"968\u00d74=3872\"],\n  [\"730\u00d78=5840\", \"648\u00d73=1944\"],\n  [\"857\u00d78=6856\", \"619\u00d72=1238\"],\n  [\"320\u00d77=2240\", \"898\u00d73=2694\"],\n  [\"807\u00d79=7263\", \"664\u00d73=1992\"],\n  [\"850\u00d75=4250\", \"577\u00d79=5193\"],\n  [\"427\u00d78=3416\", \"956\u00d78=7648\"],\n  [\"584\u00d74=2336\", \"393\u00d73=1179\"],\n  [\"529\u00d76=3174\", \"429\u00d78=3432\"],\n  [\"760\u00d78=6080\", \"639\u00d72=1278\"],\n  [\"935\u00d79=8415\", \"421\u00d72=842\"],\n  [\"416\u00d76=2496\", \"307\u00d77=2149\"],\n  [\"568\u00d76=3408\", \"467\u00d79=4203\"],\n  [\"731\u00d73=2193\", \"972\u00d75=4860\"],\n  [\"259\u00d77=1813\", \"540\u00d73=1620\"],\n  [\"599\u00d78=4792\", \"760\u00d73=2280\"],\n  [\"241\u00d72=482\", \"578\u00d77=4046\"],\n  [\"559\u00d73=1677\", \"891\u00d78=7128\"],\n  [\"175\u00d78=1400\", \"822\u00d75=4110\"],\n  [\"498\u00d77=3486\", \"786\u00d75=3930\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const rng of results.items) {\n    rng.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each \"A\u00d7B=C\" table-cell expression with its\n# updated value, per the commit diff. wdReplaceAll (2) is safe here because\n# every old value is a unique substring in the document.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-22 Monday\", \"2025-09-23 Tuesday\"),\n    @(\"516\u00d77=3612\", \"716\u00d73=2148\"),\n    @(\"425\u00d72=850\", \"467\u00d75=2335\"),\n    @(\"838\u00d73=2514\", \"135\u00d73=405\"),\n    @(\"824\u00d74=3296\", \"715\u00d73=2145\"),\n    @(\"531\u00d77=3717\", \"245\u00d79=2205\"),\n    @(\"154\u00d77=1078\", \"968\u00d74=3872\"),\n    @(\"730\u00d78=5840\", \"648\u00d73=1944\"),\n    @(\"857\u00d78=6856\", \"619\u00d72=1238\"),\n    @(\"320\u00d77=2240\", \"898\u00d73=2694\"),\n    @(\"807\u00d79=7263\", \"664\u00d73=1992\"),\n    @(\"850\u00d75=4250\", \"577\u00d79=5193\"),\n    @(\"427\u00d78=3416\", \"956\u00d78=7648\"),\n    @(\"584\u00d74=2336\", \"393\u00d73=1179\"),\n    @(\"529\u00d76=3174\", \"429\u00d78=3432\"),\n    @(\"760\u00d78=6080\", \"639\u00d72=1278\"),\n    @(\"935\u00d79=8415\", \"421\u00d72=842\"),\n    @(\"416\u00d76=2496\", \"307\u00d77=2149\"),\n    @(\"568\u00d76=3408\", \"467\u00d79=4203\"),\n    @(\"731\u00d73=2193\", \"972\u00d75=4860\"),\n    @(\"259\u00d77=1813\", \"540\u00d73=1620\"),\n    @(\"599\u00d78=4792\", \"760\u00d73=2280\"),\n    @(\"241\u00d72=482\", \"578\u00d77=4046\"),\n    @(\"559\u00d73=1677\", \"891\u00d78=7128\"),\n    @(\"175\u00d78=1400\", \"822\u00d75=4110\"),\n    @(\"498\u00d77=3486\", \"786\u00d75=3930\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $result = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        throw \"Replacement failed for: $($pair[0])\"\n    }\n}\n"}
